$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first 3 data rows (rows 2-4), shifting remaining rows up.
$ws.Rows("2:4").Delete()
